# 성능 분석표.xlsx - add "Merge-Ins" benchmark column (G) next to "Merge" (F)
# and restore the MergeSortVisualizer data, per commit:
#   "MergeSortVisualizer 사용으로 원상복귀 / Excel 업데이트"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone F1:F23's formatting (style s="1") into the new G1:G23 range so the
# new column visually matches the rest of the table (same cell style index).
$ws.Range("F1:F23").Copy()
$ws.Range("G1:G23").PasteSpecial(-4122)  # xlPasteFormats

# Header
$ws.Range("G1").Value = "Merge-Ins"

# Benchmark timings (seconds) for the Merge-Ins sort, one per input size row.
$gValues = @(
    0,        # row 2  - 100
    0.002,    # row 3  - 1000
    0.004,    # row 4  - 2000
    0.004,    # row 5  - 3000
    0.005,    # row 6  - 4000
    0.007,    # row 7  - 5000
    0.008,    # row 8  - 6000
    0.012,    # row 9  - 7000
    0.012,    # row 10 - 8000
    0.017,    # row 11 - 9000
    0.019,    # row 12 - 10000
    0.029,    # row 13 - 15000
    0.041,    # row 14 - 20000
    0.069,    # row 15 - 30000
    0.097,    # row 16 - 40000
    0.132,    # row 17 - 50000
    0.212,    # row 18 - 100000
    0.564,    # row 19 - 200000
    0.801,    # row 20 - 300000
    1.25,     # row 21 - 400000
    1.699,    # row 22 - 500000
    3.231     # row 23 - 1000000
)

for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}

# Row 24 (10000000) has no Merge-Ins measurement, same as it has no Count (C24) value.

# Restore the selection to where the author left off.
[void]$ws.Range("N13").Select()
